# Automatische test-sync: 2025-06-27 22:49:50
# Adds a new "Testmail #4" log row to the Logs sheet, extends the
# conditional-formatting ranges to include it, and refreshes the
# Dashboard category summary (counts + sort order) to account for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- 1. Append the new log entry as row 13 -----------------------------
$ws.Range("A13").Value = "Wil je dit artikel voor me inkopen?"
$ws.Range("B13").Value = "mailmind.test@zohomail.eu"
$ws.Range("C13").Value = "Testmail #4: Wil je dit artikel voor me inkopen?"
$ws.Range("D13").Value = "Bestelling / Levering"
$ws.Range("E13").Value = "Beste klant,`nBedankt voor je interesse in ons artikel. Helaas kan ik je op basis van dit bericht niet verder helpen. Kun je meer details geven over welk artikel je wilt inkopen en op welke manier? Zo kan ik je beter assisteren.`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent"
$ws.Range("F13").Value = "2025-06-27 22:49:08"
$ws.Range("G13").Value = "Ja"
$ws.Range("H13").Value = "Ja"
$ws.Range("I13").Value = "Nee"

# The multi-line answer in E13 triggers an implicit custom row height;
# auto-fit it back down so row 13 carries no explicit height override,
# matching the other data rows.
$ws.Rows.Item(13).AutoFit()

# --- 2. Extend the conditional formatting ranges to row 13 -------------
$ws.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D13"))
$ws.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G13"))
$ws.Range("H2:H12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H13"))
$ws.Range("I2:I12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I13"))

# --- 3. Refresh the Dashboard category roll-up --------------------------
# "Bestelling / Levering" now has 2 entries and moves up (sorted by
# descending count), pushing "Planning / Afspraak" and "Overig" down one
# row each (their counts stay at 1).
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A3").Value = "Bestelling / Levering"
$dash.Range("B3").Value = 2
$dash.Range("A4").Value = "Planning / Afspraak"
$dash.Range("B4").Value = 1
$dash.Range("A5").Value = "Overig"
$dash.Range("B5").Value = 1
